# Update version/build strings from the January 30 nightly build to the
# version 1.0.0 (Feb 3 2026) release build.

$oldVersionString = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersionString = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: header line in A2
$wsAbout.Range("A2").Value = "Version: " + $newVersionString

# "About" sheet: citation text in A6
$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Dendrobium Coal Mine, Australia, M0034, version ''' + $newVersionString + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$wsAbout.Range("A6").Value = $newCitation

# "Boundaries and methane sources" sheet: build_version column (S) for rows 2-8
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsBoundaries.Cells.Item($row, 19)
    if ($cell.Value() -eq $oldVersionString) {
        $cell.Value = $newVersionString
    }
}
